$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-11-11 Monday"; new = "2024-11-12 Tuesday"},
    @{old = "15×51="; new = "94×91="},
    @{old = "50×95="; new = "67×42="},
    @{old = "76×17="; new = "82×67="},
    @{old = "78×27="; new = "89×31="},
    @{old = "32×13="; new = "50×15="},
    @{old = "76×15="; new = "64×55="},
    @{old = "42×36="; new = "91×21="},
    @{old = "20×80="; new = "41×39="},
    @{old = "65×66="; new = "99×51="},
    @{old = "43×69="; new = "91×56="},
    @{old = "33×91="; new = "20×63="},
    @{old = "56×57="; new = "27×63="},
    @{old = "28×80="; new = "62×82="},
    @{old = "49×63="; new = "96×94="},
    @{old = "51×20="; new = "46×16="},
    @{old = "81×69="; new = "72×29="},
    @{old = "82×56="; new = "46×96="},
    @{old = "18×60="; new = "79×61="},
    @{old = "78×96="; new = "92×69="},
    @{old = "64×48="; new = "51×40="},
    @{old = "81×59="; new = "49×71="},
    @{old = "22×60="; new = "21×40="},
    @{old = "72×44="; new = "77×43="},
    @{old = "84×24="; new = "18×20="},
    @{old = "33×94="; new = "47×94="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
